$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Add new classification rows for AX_Bauwerksfunktion_BauwerkOderAnlage-
# FuerSportFreizeitUndErholung: building-function values 1431, 1432 and
# 1650 are missing from the table and need to be inserted.
# (1441 / 1442 already exist.)
# -----------------------------------------------------------------------

# 1) Make room for the two new "Zuschauertribuene" rows right after the
#    existing "Zuschauertribuene" row (row 2). Copy row 2 twice so the
#    new rows inherit the same cell formatting (A = wrap style, C =
#    hyperlink style) as the rest of the table.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# 2) Fill in the two new rows (3 and 4)
$ws.Range("A3").Value = "Zuschauertribüne, überdacht"
$ws.Range("B3").Value = 1431
$ws.Range("B3").ClearFormats()

$ws.Range("A4").Value = "Zuschauertribüne, nicht überdacht"
$ws.Range("B4").Value = 1432
$ws.Range("B4").ClearFormats()

# (C3 / C4 already contain the correct "stadium" codelist URL plus
#  hyperlink-style formatting, copied along with the rest of row 2.)

# 3) The former row 6 (Sprungschanze (Anlauf), value 1470) is now row 8
#    after the insert above. Its "Primär" cell (C8) used the plain
#    style instead of the hyperlink style used everywhere else in the
#    column - align it with the rest of the column.
$ws.Range("C5").Copy()
$ws.Range("C8").PasteSpecial(-4122) | Out-Null

# 4) Former (empty) row 8 is now row 10. Turn it into the new
#    "Wassersportanlage" entry (value 1650) - its existing formatting
#    already matches the target styling, so only the values need to be
#    set.
$ws.Range("A10").Value = "Wassersportanlage"
$ws.Range("B10").Value = 1650
$ws.Range("C10").Value = "{{project:BUILDINGNATURE}}sonstiges"

# 5) The two now-empty filler rows 11 and 12 (previously rows 9 and 10)
#    are no longer needed - remove them so the rest of the sheet keeps
#    its original row numbers.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(11).Delete()

$excel.CutCopyMode = 0
